$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 27
$ws1.Range("F4").Value = 403
$ws1.Range("F5").Value = 5115
$ws1.Range("F6").Value = 5115
$ws1.Range("F7").Value = 106
$ws1.Range("F11").Value = 1155
$ws1.Range("F12").Value = 707
$ws1.Range("F13").Value = 4966
$ws1.Range("F15").Value = 63
$ws1.Range("F17").Value = 210
$ws1.Range("I17").Value = "//i1.hdslb.com/bfs/openplatform/202406/mJt8McPp1718594709773.jpeg"
$ws1.Range("F18").Value = 223
$ws1.Range("F21").Value = 3769
$ws1.Range("F24").Value = 3657
$ws1.Range("F25").Value = 174
$ws1.Range("F30").Value = 202
$ws1.Range("F36").Value = 6485
$ws1.Range("F37").Value = 1031
$ws1.Range("F38").Value = 486
$ws1.Range("F42").Value = 1325
$ws1.Range("F44").Value = 646
$ws1.Range("F46").Value = 2223
$ws1.Range("F49").Value = 764
$ws1.Range("F50").Value = 908

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 132
$ws2.Range("F9").Value = 81
$ws2.Range("F23").Value = 802

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 27
$ws4.Range("F7").Value = 403
$ws4.Range("F8").Value = 5115
$ws4.Range("F9").Value = 5115
$ws4.Range("F10").Value = 106
$ws4.Range("F13").Value = 81
$ws4.Range("F15").Value = 1155
$ws4.Range("F16").Value = 707
$ws4.Range("F17").Value = 4966
$ws4.Range("F19").Value = 63
$ws4.Range("F21").Value = 210
$ws4.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202406/mJt8McPp1718594709773.jpeg"
$ws4.Range("F22").Value = 223
$ws4.Range("F25").Value = 3769
$ws4.Range("F26").Value = 3657
$ws4.Range("F27").Value = 174
$ws4.Range("F31").Value = 202
$ws4.Range("F37").Value = 6485
$ws4.Range("F38").Value = 1031
$ws4.Range("F39").Value = 486
$ws4.Range("F43").Value = 1325
$ws4.Range("F45").Value = 646
$ws4.Range("F46").Value = 2223
$ws4.Range("F48").Value = 764
$ws4.Range("F49").Value = 908
